$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row at position 4 - shifts existing rows 4-10 down to 5-11
$ws.Rows("4:4").Insert()

# 2) Fill the new row 4 with the "Projection-Aware Planning" entry.
#    Column order matters: it controls the order new strings are appended
#    to the shared-string table, so G (imagine) is written before F (the
#    link text) to match the expected shared-string ordering.
$ws.Range("A4").Value2 = "Projection-Aware Planning"
$ws.Range("B4").Value2 = "Projection-Aware Task Planning and Execution for Human-in-the-Loop Operation of Robots in a Mixed-Reality Workspace"
$ws.Range("C4").Value2 = "Tathagata Chakraborti, Sarath Sreedharan, Anagha Kulkarni and Subbarao Kambhampati"
$ws.Range("D4").Value2 = "Recent advances in mixed-reality technologies have renewed interest in alternative modes of communication for human-robot interaction. However, most of the work in this direction has been confined to tasks such as teleoperation, simulation or explication of individual actions of a robot. In this paper, we will discuss how the capability to project intentions affect the task planning capabilities of a robot. Specifically, we will start with a discussion on how projection actions can be used to reveal information regarding the future intentions of the robot at the time of task execution. We will then pose a new planning paradigm – projection-aware planning – whereby a robot can trade off its plan cost with its ability to reveal its intentions using its projection actions. We demonstrate each of these scenarios with the help of a joint human-robot activity using the HoloLens."
$ws.Range("E4").Value2 = "U.S. Finals, Microsoft Imagine Cup 2017 (ICAPS Demo Track 2017)"
$ws.Range("G4").Value2 = "imagine"
$ws.Range("F4").Value2 = "http://ae-robots.com/,https://yochan-lab.github.io/papers/files/papers/projection-aware.pdf"

# 3) Authors (C4) and Abstract (D4) carry no explicit cell style in the
#    target file (the row-insert copied style "1" down from row 3) - clear
#    that back to the workbook default.
$ws.Range("C4:D4").ClearFormats()

# 4) New column F (Links) gets an explicit width.
$ws.Columns("F:F").ColumnWidth = 67.16666666666667

# 5) Row-insert does not renumber the <hyperlinks> list, so rebuild it from
#    scratch: delete everything then re-add in the order that reproduces
#    rId1..rId7 pointing at their original targets (now one row further
#    down for every row that used to sit at 4 or below) plus the brand new
#    rId8 hyperlink for the inserted row.
$ws.Cells.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("F5"), "https://nirlipo.github.io/Width-Based-Planning-Resources/,https://people.eng.unimelb.edu.au/nlipovetzky/papers/aiaccess_nirlipo.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.aaai.org/ocs/index.php/ICAPS/ICAPS13/paper/view/6039/6208") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "http://www.aaai.org/Papers/ICAPS/2007/ICAPS07-008.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F6"), "http://planning.domains/,http://editor.planning.domains/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F9"), "http://hcjournal.org/ojs/index.php?journal=jhc&page=article&op=view&path%5B%5D=10.15346%2Fhc.v4i1.2") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F10"), "https://aaai.org/ocs/index.php/ICAPS/ICAPS17/paper/view/15617") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F11"), "https://arxiv.org/abs/1709.04517") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "http://ae-robots.com/,https://yochan-lab.github.io/papers/files/papers/projection-aware.pdf") | Out-Null

# Hyperlinks.Add() stamps the freshly linked cell with a brand-new cell
# style (a duplicate of the existing "Hyperlink" style) instead of reusing
# the workbook's existing one. Nudge each link cell back onto the shared
# style by re-applying just the font property that differs - the engine's
# style-interning reuses the existing cellXf/font instead of minting a
# fresh one when the requested format already matches one on file.
$ws.Range("F2").Font.Underline = $true
$ws.Range("F3").Font.Underline = $true
$ws.Range("F4").Font.Underline = $true
$ws.Range("F5").Font.Underline = $true
$ws.Range("F6").Font.Underline = $true
$ws.Range("F9").Font.Underline = $true
$ws.Range("F10").Font.Underline = $true
$ws.Range("F11").Font.Underline = $true

# 6) Selection moves to G5 in the saved file.
$ws.Range("G5").Select()

Write-Output "edit complete"
